# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (AD1:AF1) -- copy formatting (bold, centered, bordered)
# from the existing header style used on row 1 (e.g. AC1), then set text.
$ws.Range("AD1:AF1").Value = "Wins"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every data row (2-54): Wins=85, Losses=77, Ties=0
$lastRow = 54
$ws.Range("AD2:AD$lastRow").Value = 85
$ws.Range("AE2:AE$lastRow").Value = 77
$ws.Range("AF2:AF$lastRow").Value = 0
